# Update the 'timestamp' column (Z) on the Log_Muestras sheet with the new
# run timestamps recorded when this dataset augmentation log was
# regenerated (dataset Us Crime agregado).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$ws.Range("Z2:Z7").Value = "2025-11-13T06:52:05.067726"
$ws.Range("Z8:Z9").Value = "2025-11-13T06:52:05.068817"
$ws.Range("Z10:Z12").Value = "2025-11-13T06:52:05.069274"
$ws.Range("Z13").Value = "2025-11-13T06:52:05.069685"
$ws.Range("Z14:Z57").Value = "2025-11-13T06:52:05.069748"
$ws.Range("Z58:Z69").Value = "2025-11-13T06:52:05.166490"
$ws.Range("Z70").Value = "2025-11-13T06:52:05.179692"
$ws.Range("Z71:Z73").Value = "2025-11-13T06:52:05.331627"
$ws.Range("Z74:Z79").Value = "2025-11-13T06:52:05.332625"
